$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a brand-new first paragraph with the YouTube reference
# link, bold + Segoe UI Emoji, "Reff" wrapped in spell-check proof markers.
# ---------------------------------------------------------------------------
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$refXml = '<w:p ' + $xmlNs + '>' +
    '<w:pPr><w:rPr>' +
      '<w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>' +
      '<w:b/><w:bCs/>' +
    '</w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' +
      '<w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>' +
      '<w:b/><w:bCs/>' +
    '</w:rPr><w:t>Reff</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' +
      '<w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>' +
      '<w:b/><w:bCs/>' +
    '</w:rPr><w:t>: https://www.youtube.com/watch?v=ilgpzlE7Hds&amp;list=PLkGgF5-rBAbi4c4LbVExbh62S_kcGjVHb</w:t></w:r>' +
  '</w:p>'

$startRange = $d.Range(0, 0)
$startRange.InsertXML($refXml)

Write-Output "done"
